$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the interpolation input parameters (xmin/xmax row changes)
$ws.Range("B2").Value = 35
$ws.Range("B3").Value = 400

# Move the active selection to B2 (matches the saved selection in the diff)
$ws.Range("B2").Select()
